# Apply changes described by the commit "add tpe.suggest parameters to config file"
# on the st-architecture.pptx slide (Slide 2 of the presentation):
#  1. Remove the "Elbow Connector 1050" connector shape.
#  2. Remove the "TextBox 138" shape (the SmartTuning configmap reload note).
#  3. Add an end-arrowhead (triangle) line to the "Freeform 54" shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# 1. Delete the elbow connector shape.
$elbow = $s.Shapes.Item("Elbow Connector 1050")
$elbow.Delete()

# 2. Delete the now-orphaned "TextBox 138" explanatory text box.
$textBox138 = $s.Shapes.Item("TextBox 138")
$textBox138.Delete()

# 3. Give the "Freeform 54" freeform shape a visible line with a triangle
#    arrowhead on the tail end (no arrowhead on the head end), both using
#    medium width/length, matching the target OOXML <a:ln> block.
$freeform54 = $s.Shapes.Item("Freeform 54")
$freeform54.Line.Visible = $true
$freeform54.Line.BeginArrowheadStyle = 1
$freeform54.Line.BeginArrowheadWidth = 2
$freeform54.Line.BeginArrowheadLength = 2
$freeform54.Line.EndArrowheadStyle = 2
$freeform54.Line.EndArrowheadWidth = 2
$freeform54.Line.EndArrowheadLength = 2
